# Update Angola MSME "Enterprises density (per 1000 people)" row:
# Micro column (B13): "1" -> "0.96"
# SMEs column  (C13): "0" -> "0.04"
#
# These values need to stay stored as text (shared-string) cells, matching
# how they were originally authored, rather than being auto-coerced into
# numeric cells by a plain Value assignment. We build the numeric-looking
# text in a scratch cell via a text formula, then copy/paste-special just
# the values back onto the target cells (mirrors the real-Excel trick of
# pasting computed text without dragging the formula/format along), and
# finally clear the scratch cell so it leaves no trace in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Formula = '="0.96"'
$scratch.Copy()
$ws.Range("B13").PasteSpecial(-4163)  # xlPasteValues

$scratch.Formula = '="0.04"'
$scratch.Copy()
$ws.Range("C13").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
